$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at the very top of the sheet. This pushes the
# existing header row (old row 1) down to row 3, and every data row down by
# two (old row 2 -> row 4, old row 3 -> row 5, ... old row 36 -> row 38).
$ws.Range("A1:N2").Insert()

# Row 3 now holds the old header text but still carries the old row 1's
# bold / centered / bordered header style. Copy that style up onto the new
# row 1 (which will become the numeric index header), then strip it from
# row 3 so it becomes a plain, unstyled data row like the other rows.
$ws.Range("A3:N3").Copy()
$ws.Range("A1:N1").PasteSpecial(-4122)
$ws.Range("A3:N3").ClearFormats()
$excel.CutCopyMode = 0

# New row 1: a numeric column-index header, 0 through 13.
$ws.Range("A1").Value = 0
$ws.Range("B1").Value = 1
$ws.Range("C1").Value = 2
$ws.Range("D1").Value = 3
$ws.Range("E1").Value = 4
$ws.Range("F1").Value = 5
$ws.Range("G1").Value = 6
$ws.Range("H1").Value = 7
$ws.Range("I1").Value = 8
$ws.Range("J1").Value = 9
$ws.Range("K1").Value = 10
$ws.Range("L1").Value = 11
$ws.Range("M1").Value = 12
$ws.Range("N1").Value = 13

# New row 2: blank across the board except E2, which reads "Washer".
$ws.Range("A2").Value = ""
$ws.Range("B2").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("D2").Value = ""
$ws.Range("E2").Value = "Washer"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = ""
$ws.Range("I2").Value = ""
$ws.Range("J2").Value = ""
$ws.Range("K2").Value = ""
$ws.Range("L2").Value = ""
$ws.Range("M2").Value = ""
$ws.Range("N2").Value = ""
